# feat: add 2022-Q1 data
#
# Before: sheets are 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计 (summary table:
#   日期 / 持有数量(只) / 持有市值(亿元), one row per quarter).
#
# After: a new "2022-Q1" sheet (holdings detail, same layout as the other
#   quarter sheets) is inserted right before "总计", and "总计" gains a new
#   top row summarizing the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. The existing "总计" sheet becomes "2022-Q1" (keeps its sheetId/part),
#    and a brand-new "总计" sheet is inserted right after it to hold the
#    (updated) summary table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------
# 2. Build the new "总计" summary table: same 3 columns as before, plus a
#    new first data row for 2022-Q1, with the older rows pushed down.
# ---------------------------------------------------------------------
$totalHeader = $wb.Worksheets.Item("2021-Q4").Range("B1")
$totalHeader.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("A2:A6").PasteSpecial(-4122)   # reuse same index-column format

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @(0, "2022-Q1", 6, 0.99),
    @(1, "2021-Q4", 2, 1.4),
    @(2, "2021-Q3", 5, 0.11),
    @(3, "2021-Q2", 2, 0.05),
    @(4, "2021-Q1", 3, 0.09)
)
$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Rebuild the "2022-Q1" sheet itself as a holdings-detail sheet (same
#    column layout as the other quarter sheets): clear the old summary
#    content first, pick up the header/body formatting from "2021-Q4"
#    (which already uses the "基金规模" wording we need), then fill in
#    the fund-holding rows for 2022-Q1.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")

$q1.Cells.Clear()

$template.Range("A1:H1").Copy()
$q1.Range("A1:H1").PasteSpecial(-4122)   # xlPasteFormats (header row)

$template.Range("A2:H2").Copy()
$q1.Range("A2:H7").PasteSpecial(-4122)   # xlPasteFormats, tiles cleanly (1-row source)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$holdings = @(
    @(0, "014207", "华安产业精选混合A",     "27.31", "62.03", "1.67", "0.4561", 10),
    @(1, "159611", "广发中证全指电力ETF",     "13.38", "99.14", "2.30", "0.3077", 9),
    @(2, "014208", "华安产业精选混合C",     "7.93",  "62.03", "1.67", "0.1324", 10),
    @(3, "080005", "长盛量化红利混合",       "2.66",  "69.88", "2.69", "0.0716", 5),
    @(4, "008778", "嘉实中证500指数增强A", "0.93",  "93.42", "1.89", "0.0176", 9),
    @(5, "008779", "嘉实中证500指数增强C", "0.35",  "93.42", "1.89", "0.0066", 9)
)

$r = 2
foreach ($row in $holdings) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
